# Journal de bord.xlsx - "Correction de la structure des dossiers"
#
# 1. Add a new journal entry (13.03.2018 - "Mise en place du serveur Web")
#    as row 11 of Brian's "Log_5" table, extending B2:C10 -> B2:C11.
# 2. Switch the active/selected tab from "Alexandre" to "Brian" and update
#    the on-sheet selections accordingly.

$wb = $excel.ActiveWorkbook

$wsAlexandre = $wb.Worksheets.Item("Alexandre")
$wsBrian     = $wb.Worksheets.Item("Brian")

# --- 1. Append the new row to Brian's journal table ------------------------

$table = $wsBrian.ListObjects.Item("Log_5")
$table.ListRows.Add() | Out-Null

$wsBrian.Cells.Item(11, 2).Value = 43172
$wsBrian.Cells.Item(11, 3).Value = "Mise en place du serveur Web"

# Copy the formatting (styles) of the row above onto the freshly added row
# so the new cells keep the "Date Column" / "Event Column" look instead of
# picking up the default style.
$wsBrian.Range("B10:C10").Copy()
$wsBrian.Range("B11:C11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsBrian.Rows.Item(11).RowHeight = $wsBrian.Rows.Item(10).RowHeight

# --- 2. Move the active tab from Alexandre to Brian -------------------------

# Alexandre scrolls down a bit (top-left visible cell becomes A4) and is no
# longer the selected tab.
$wsAlexandre.Activate()
$wsAlexandre.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$wsAlexandre.Range("C12").Select()

# Brian becomes the active/selected sheet, with the selection sitting on the
# next empty row right below the new entry.
$wsBrian.Activate()
$wsBrian.Range("C12").Select()
